# Happy Farm Fruit liquidation workbook - pallet numbers correction.
#
# The "pallet NO." column (D10:D17) held truncated 4-digit numbers
# (e.g. 1817) that Excel was treating as quantitative data. They are
# actually identifiers, so they are rewritten as their real 7-digit
# pallet numbers (151xxxx) and stored as text, not numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> correct pallet number (stored as text).
$palletNumbers = [ordered]@{
    "D10" = "1511817"
    "D11" = "1511818"
    "D12" = "1511816"
    "D13" = "1511817"
    "D14" = "1511818"
    "D15" = "1511819"
    "D16" = "1511861"
    "D17" = "1511816"
}

foreach ($addr in $palletNumbers.Keys) {
    $cell = $ws.Range($addr)

    # Remember the cell's current (numeric) display format so we can
    # restore it after the write - only the underlying type should
    # change from number to text, not the visual format.
    $originalFormat = $cell.NumberFormat

    # Forcing a text format while writing keeps the literal digits
    # intact (no scientific notation / leading-zero loss) and makes
    # Excel store the value as a real string instead of re-parsing it
    # back into a number.
    $cell.NumberFormat = "@"
    $cell.Value = $palletNumbers[$addr]
    $cell.NumberFormat = $originalFormat
}

# Leave the selection where the edit left off, as in the source session.
$ws.Range("D20").Select()
